# Applies the cryptos.xlsx price/volume refresh described in the commit
# "Updated cryptos list on Sat Sep 23 07:52:03 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.665.85"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "1.598.97"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0619"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.59"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.33%  "
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("D12").Value = "1.823.08"
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("D13").Value = "1.585.28"
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.94"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.31%  "
$ws.Range("D17").Value = "26.654.34"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("D18").Value = "0.0₃0734"
$ws.Range("E18").Value = "  +0.56%  "
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("E20").Value = "  -0.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.11%  "
$ws.Range("E23").Value = "  +0.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.99%  "
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.23%  "
$ws.Range("E28").Value = "  -0.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.31"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  +1.89%  "
$ws.Range("E31").Value = "  +0.16%  "
$ws.Range("E32").Value = "  +0.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.46%  "
$ws.Range("D34").Value = "1.279.82"
$ws.Range("E34").Value = "  -1.11%  "
$ws.Range("E35").Value = "  -10.07%  "
$ws.Range("E36").Value = "  +0.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.49"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0170"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.50%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.835"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.50%  "
$ws.Range("B40").Value = "WEMIXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.06"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +17.63%  "
$ws.Range("E41").Value = "  +2.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.785"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.60%  "
$ws.Range("D45").Value = "1.735.41"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.57"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.91%  "
$ws.Range("E47").Value = "  -2.51%  "
$ws.Range("E48").Value = "  +3.10%  "
$ws.Range("E49").Value = "  +1.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.42"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.02%  "
